$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New URL referenced by the "Alianza del Pacífico" row (row 6)
$url = "https://alianzapacifico.net/download/decision-no-9-anexo-suplementos-alimenticios-eliminacion-de-obstaculos-tecnicos/"

# Put the URL text in E6 and turn it into a working hyperlink
$ws.Range("E6").Value = $url
$ws.Hyperlinks.Add($ws.Range("E6"), $url)

# Match the "Hipervínculo" style already used by the other link cells (E2, E3, E5)
$ws.Range("E6").Style = $ws.Range("E5").Style

# Update the visible selection to E7 (row just below the new data)
[void]$ws.Range("E7").Select()
